# Notice of Entry template: make the title a real Heading 1 and replace
# every static label with its {{PLACEHOLDER}} counterpart.
#
# Each paragraph is rewritten via Range.InsertXML with a minimal
# WordprocessingML fragment (wrapped in the flat-OPC pkg:package envelope
# that Word's COM InsertXML expects). Driving the exact paragraph XML this
# way - rather than mutating Font.Bold/Style on the existing run, or doing a
# plain Find/Replace of the visible text - avoids leaving stray artifacts
# behind, such as:
#   - an empty <w:rPr/> left on the run after flipping paragraph style
#   - the old trailing <w:br/> run surviving a text-only replace
#
# One wrinkle: InsertXML-ing the *entire* Range of the document's last
# paragraph (whose Range.End coincides with Content.End) leaves a spurious
# trailing empty paragraph behind. Excluding the final paragraph-mark
# character from the replaced range avoids that, but it also means any
# <w:pPr> in the inserted XML is ignored (paragraph formatting lives on the
# paragraph mark). Since only the first paragraph needs a <w:pPr> change,
# and it isn't the last paragraph, that's not a problem here.

$d = $word.ActiveDocument

function New-ParagraphXml([string]$innerPPr, [string]$text) {
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerPPr + '<w:r><w:t>' + $escaped + '</w:t></w:r></w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphText([int]$index, [string]$text, [string]$innerPPr = "") {
    $p = $d.Paragraphs.Item($index)
    $range = $p.Range
    # Trim the trailing paragraph-mark character off the very last paragraph
    # in the story so InsertXML replaces content in place instead of
    # appending a new (empty) paragraph after it.
    if ($range.End -eq $d.Content.End) {
        $range = $d.Range($range.Start, $range.End - 1)
    }
    [void]$range.InsertXML((New-ParagraphXml $innerPPr $text))
}

# 1. Title -> real Heading 1 style, plain run (drops the manual b/sz=32 run formatting).
Set-ParagraphText 1 "Notice of Entry" '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>'

# 2. Date: 2026-01-21 -> Tenant Name placeholder
Set-ParagraphText 2 "Tenant Name: {{TENANT_NAME}}"

# 3. "Not legal advice..." (+ trailing <w:br/>) -> Property Address placeholder, no break
Set-ParagraphText 3 "Property Address: {{PROPERTY_ADDRESS}}"

# 4. Tenant: -> Unit placeholder
Set-ParagraphText 4 "Unit: {{UNIT_NUMBER}}"

# 5. Property: -> Date of Notice placeholder
Set-ParagraphText 5 "Date of Notice: {{NOTICE_DATE}}"

# 6. Date of Entry: -> Planned Entry Date/Time placeholder
Set-ParagraphText 6 "Planned Entry Date/Time: {{ENTRY_DATE_TIME}}"

# 7. Reason: -> Reason for Entry placeholder
Set-ParagraphText 7 "Reason for Entry: {{REASON_FOR_ENTRY}}"

# 8. Contact: -> Landlord/Manager placeholder (last paragraph in the story)
Set-ParagraphText 8 "Landlord/Manager: {{LANDLORD_NAME}}"

Write-Host "Notice of Entry template updated."
